$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 25217
$ws.Range("D2").Value = 1648
$ws.Range("B3").Value = 23920
$ws.Range("D3").Value = 1761
$ws.Range("B4").Value = 23123
$ws.Range("D4").Value = 1929
$ws.Range("B5").Value = 22966
$ws.Range("D5").Value = 2031.5
$ws.Range("B6").Value = 23289
$ws.Range("D6").Value = 2022
$ws.Range("B7").Value = 24770
$ws.Range("D7").Value = 2026
$ws.Range("B8").Value = 29927
$ws.Range("C8").Value = 229
$ws.Range("D8").Value = 1989
$ws.Range("B9").Value = 34860
$ws.Range("C9").Value = 2139
$ws.Range("D9").Value = 1826.5
$ws.Range("B10").Value = 38009
$ws.Range("C10").Value = 5590
$ws.Range("D10").Value = 1683
$ws.Range("B11").Value = 38538
$ws.Range("C11").Value = 8582
$ws.Range("D11").Value = 1703.5
$ws.Range("B12").Value = 38173
$ws.Range("C12").Value = 10384
$ws.Range("D12").Value = 1862.5
$ws.Range("B13").Value = 37920
$ws.Range("C13").Value = 10970
$ws.Range("D13").Value = 2063.5
$ws.Range("B14").Value = 36343
$ws.Range("C14").Value = 10639
$ws.Range("D14").Value = 2264
$ws.Range("B15").Value = 36267
$ws.Range("C15").Value = 9289
$ws.Range("D15").Value = 2452
$ws.Range("B16").Value = 37320
$ws.Range("C16").Value = 6815
$ws.Range("D16").Value = 2520
$ws.Range("B17").Value = 37901
$ws.Range("C17").Value = 3576
$ws.Range("D17").Value = 2488
$ws.Range("B18").Value = 38499
$ws.Range("C18").Value = 898
$ws.Range("D18").Value = 2456.5
$ws.Range("B19").Value = 39551
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 2439
$ws.Range("B20").Value = 39337
$ws.Range("D20").Value = 2454
$ws.Range("B21").Value = 39005
$ws.Range("D21").Value = 2517.5
$ws.Range("B22").Value = 36911
$ws.Range("D22").Value = 2628
$ws.Range("B23").Value = 33958
$ws.Range("D23").Value = 2822
$ws.Range("B24").Value = 30490
$ws.Range("D24").Value = 3039.5
$ws.Range("B25").Value = 27615
$ws.Range("D25").Value = 3276.5
